# Update "想去人数" (F column) counts on sheets 展览 (sheet1), 演出 (sheet2),
# and 全部类型 (sheet4) to reflect the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheets index 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 9898
    4  = 2517
    6  = 277
    8  = 479
    9  = 726
    11 = 1219
    12 = 1012
    13 = 3061
    14 = 2308
    16 = 2005
    17 = 245
    19 = 484
    21 = 529
    22 = 42
    23 = 209
    26 = 356
    28 = 339
    29 = 546
    31 = 193
    32 = 1557
    33 = 257
    34 = 1587
    35 = 79
    36 = 385
    37 = 39
    38 = 415
    39 = 864
    41 = 335
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (Worksheets index 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    2 = 27
    5 = 7
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" (Worksheets index 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 9898
    4  = 2517
    5  = 27
    8  = 277
    10 = 479
    11 = 726
    13 = 1219
    14 = 1012
    15 = 3061
    16 = 2308
    18 = 2005
    19 = 245
    21 = 484
    23 = 529
    24 = 42
    25 = 209
    28 = 356
    30 = 339
    31 = 546
    33 = 7
    36 = 193
    37 = 1557
    39 = 257
    40 = 1587
    41 = 79
    43 = 385
    44 = 39
    45 = 415
    46 = 864
    48 = 335
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
